# Update the "Förändrad" (Changed) date column (C) from 2023-09-06 (45175)
# to 2023-09-08 (45177) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet (data starts at row 2; row 1 is the header).
$lastRow = $ws.UsedRange.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 261 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C ("Förändrad")
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
